# Generate Report for Handback
# Fills in "Latest Target File" (F) and "Latest Handback File" (G) columns
# on the zh-cn and de-de sheets, updates the handback status text on the
# Overview sheet (via the shared "Status" string), and records the actual
# handback datetimes.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$urlMd1  = "https://github.com/OpenLocalizationTest/oltest/blob/1cf2ed4331334abe075ef22ca2666fc4b695498c/e2e/3db582a9-0042-49a5-a2ee-b7c3e39538de.md"
$urlMd2  = "https://github.com/OpenLocalizationTest/oltest/blob/1cf2ed4331334abe075ef22ca2666fc4b695498c/e2e/42082a3a-6b7a-442c-a062-087e4d6e6762.md"
$urlXlfZh1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b6215127f9a70ee4b31a28fee997b2e17fc86a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3db582a9-0042-49a5-a2ee-b7c3e39538de.d20dd841412fd6499ed3bd91f44bfef92db4ec25.zh-cn.xlf"
$urlXlfZh2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5b6215127f9a70ee4b31a28fee997b2e17fc86a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/42082a3a-6b7a-442c-a062-087e4d6e6762.29fa362ac2f9dcb459aaace1e688dd8ff08c2b50.zh-cn.xlf"

# Row 2 - 3db582a9-0042-49a5-a2ee-b7c3e39538de.md
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $urlMd1, [Type]::Missing, [Type]::Missing, "3db582a9-0042-49a5-a2ee-b7c3e39538de.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $urlXlfZh1, [Type]::Missing, [Type]::Missing, "3db582a9-0042-49a5-a2ee-b7c3e39538de.d20dd841412fd6499ed3bd91f44bfef92db4ec25.zh-cn.xlf")

# Row 3 - 42082a3a-6b7a-442c-a062-087e4d6e6762.md
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $urlMd2, [Type]::Missing, [Type]::Missing, "42082a3a-6b7a-442c-a062-087e4d6e6762.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $urlXlfZh2, [Type]::Missing, [Type]::Missing, "42082a3a-6b7a-442c-a062-087e4d6e6762.29fa362ac2f9dcb459aaace1e688dd8ff08c2b50.zh-cn.xlf")

# Actual handback datetime (replaces the 0001-01-01 00:00:00 placeholder)
$wsZh.Range("H2").Value = "2016-03-14 03:33:44"
$wsZh.Range("H3").Value = "2016-03-14 03:33:44"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$urlXlfDe1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08ca66ffbb76d15c6e93e04df57ec27eaaf5dbe8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3db582a9-0042-49a5-a2ee-b7c3e39538de.d20dd841412fd6499ed3bd91f44bfef92db4ec25.de-de.xlf"
$urlXlfDe2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08ca66ffbb76d15c6e93e04df57ec27eaaf5dbe8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/42082a3a-6b7a-442c-a062-087e4d6e6762.29fa362ac2f9dcb459aaace1e688dd8ff08c2b50.de-de.xlf"

# Row 2 - 3db582a9-0042-49a5-a2ee-b7c3e39538de.md
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $urlMd1, [Type]::Missing, [Type]::Missing, "3db582a9-0042-49a5-a2ee-b7c3e39538de.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $urlXlfDe1, [Type]::Missing, [Type]::Missing, "3db582a9-0042-49a5-a2ee-b7c3e39538de.d20dd841412fd6499ed3bd91f44bfef92db4ec25.de-de.xlf")

# Row 3 - 42082a3a-6b7a-442c-a062-087e4d6e6762.md
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $urlMd2, [Type]::Missing, [Type]::Missing, "42082a3a-6b7a-442c-a062-087e4d6e6762.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $urlXlfDe2, [Type]::Missing, [Type]::Missing, "42082a3a-6b7a-442c-a062-087e4d6e6762.29fa362ac2f9dcb459aaace1e688dd8ff08c2b50.de-de.xlf")

# Actual handback datetime
$wsDe.Range("H2").Value = "2016-03-14 03:33:49"
$wsDe.Range("H3").Value = "2016-03-14 03:33:49"

# ---- Overview sheet: status text (shared string used by B2/C2/B3/C3) ----
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("B2").Value = "Handed back: in sync with en-US"
$wsOv.Range("C2").Value = "Handed back: in sync with en-US"
$wsOv.Range("B3").Value = "Handed back: in sync with en-US"
$wsOv.Range("C3").Value = "Handed back: in sync with en-US"

Write-Host "Handback report generated"
